$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 8-44 with the new (shifted) field names and dtypes.
# These rows already exist in the sheet, so their formatting (e.g. the bold/
# bordered style on column A) is preserved automatically when we just set .Value.
$ws.Range("A8").Value = "homeGoals"
$ws.Range("B8").Value = "int64"
$ws.Range("A9").Value = "awayGoals"
$ws.Range("B9").Value = "int64"
$ws.Range("A10").Value = "homeGoalsHalfTime"
$ws.Range("B10").Value = "int64"
$ws.Range("A11").Value = "awayGoalsHalfTime"
$ws.Range("B11").Value = "int64"
$ws.Range("A12").Value = "home_xGoals"
$ws.Range("B12").Value = "float64"
$ws.Range("A13").Value = "home_shots"
$ws.Range("B13").Value = "int64"
$ws.Range("A14").Value = "home_shotsOnTarget"
$ws.Range("B14").Value = "int64"
$ws.Range("A15").Value = "home_deep"
$ws.Range("B15").Value = "int64"
$ws.Range("A16").Value = "home_ppda"
$ws.Range("B16").Value = "float64"
$ws.Range("A17").Value = "home_fouls"
$ws.Range("B17").Value = "int64"
$ws.Range("A18").Value = "home_corners"
$ws.Range("B18").Value = "int64"
$ws.Range("A19").Value = "home_yellowCards"
$ws.Range("B19").Value = "float64"
$ws.Range("A20").Value = "home_redCards"
$ws.Range("B20").Value = "int64"
$ws.Range("A21").Value = "home_total_assists"
$ws.Range("B21").Value = "int64"
$ws.Range("A22").Value = "home_total_xAssists"
$ws.Range("B22").Value = "float64"
$ws.Range("A23").Value = "home_total_key_passes"
$ws.Range("B23").Value = "int64"
$ws.Range("A24").Value = "home_total_xGoalsChain"
$ws.Range("B24").Value = "float64"
$ws.Range("A25").Value = "home_total_xGoalsBuildup"
$ws.Range("B25").Value = "float64"
$ws.Range("A26").Value = "home_total_yellow_cards"
$ws.Range("B26").Value = "int64"
$ws.Range("A27").Value = "home_total_red_cards"
$ws.Range("B27").Value = "int64"
$ws.Range("A28").Value = "home_total_blocked_shots"
$ws.Range("B28").Value = "float64"
$ws.Range("A29").Value = "home_total_saved_shots"
$ws.Range("B29").Value = "float64"
$ws.Range("A30").Value = "away_xGoals"
$ws.Range("B30").Value = "float64"
$ws.Range("A31").Value = "away_shots"
$ws.Range("B31").Value = "int64"
$ws.Range("A32").Value = "away_shotsOnTarget"
$ws.Range("B32").Value = "int64"
$ws.Range("A33").Value = "away_deep"
$ws.Range("B33").Value = "int64"
$ws.Range("A34").Value = "away_ppda"
$ws.Range("B34").Value = "float64"
$ws.Range("A35").Value = "away_fouls"
$ws.Range("B35").Value = "int64"
$ws.Range("A36").Value = "away_corners"
$ws.Range("B36").Value = "int64"
$ws.Range("A37").Value = "away_yellowCards"
$ws.Range("B37").Value = "float64"
$ws.Range("A38").Value = "away_redCards"
$ws.Range("B38").Value = "int64"
$ws.Range("A39").Value = "away_total_assists"
$ws.Range("B39").Value = "int64"
$ws.Range("A40").Value = "away_total_xAssists"
$ws.Range("B40").Value = "float64"
$ws.Range("A41").Value = "away_total_key_passes"
$ws.Range("B41").Value = "int64"
$ws.Range("A42").Value = "away_total_xGoalsChain"
$ws.Range("B42").Value = "float64"
$ws.Range("A43").Value = "away_total_xGoalsBuildup"
$ws.Range("B43").Value = "float64"
$ws.Range("A44").Value = "away_total_yellow_cards"
$ws.Range("B44").Value = "int64"

# Append 4 brand-new rows (45-48) for the fields that did not previously exist
# in the sheet (the list grew from 44 to 48 rows). Copy the column-A formatting
# (bold font + box border) from an existing labeled cell before setting values,
# so the new cells match the look of the rest of column A.
$ws.Range("A2").Copy()
$ws.Range("A45").PasteSpecial(-4122)
$ws.Range("A45").Value = "away_total_red_cards"
$ws.Range("B45").Value = "int64"
$ws.Range("A2").Copy()
$ws.Range("A46").PasteSpecial(-4122)
$ws.Range("A46").Value = "away_total_blocked_shots"
$ws.Range("B46").Value = "float64"
$ws.Range("A2").Copy()
$ws.Range("A47").PasteSpecial(-4122)
$ws.Range("A47").Value = "away_total_saved_shots"
$ws.Range("B47").Value = "float64"
$ws.Range("A2").Copy()
$ws.Range("A48").PasteSpecial(-4122)
$ws.Range("A48").Value = "gameresult"
$ws.Range("B48").Value = "object"

$excel.CutCopyMode = 0

Write-Host "Done updating teamstats datatype sheet."
